$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# This edit rotates the data of three data rows (24, 25, 27) — row 26
# (Knärot) is untouched in between. Effectively:
#   old row 27 (A=111957843, Rosenticka) -> new row 24
#   old row 24 (A=111958182, Järpe)      -> new row 25
#   old row 25 (A=111957798, Rosenticka) -> new row 27
# Rather than physically moving rows, we overwrite the cells of each
# row in place with the target values, and fix up the handful of
# cells that appear/disappear between the three records.
# ------------------------------------------------------------------

function Set-TextCell($range, [string]$text) {
    # Force a genuine text-typed value even when it looks numeric
    # (e.g. "3", "6"), without leaving a stray NumberFormat/style
    # behind once we're done.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

function Set-EmptyTextCell($range) {
    # Produce a present-but-empty text cell (as opposed to Clear-
    # Contents, which removes the cell entirely). Excel treats a
    # leading apostrophe as a force-text marker for an otherwise
    # blank entry; resetting the style afterwards drops the
    # quote-prefix formatting flag it introduces.
    $range.Value = "'"
    $range.Style = "Normal"
}

# ===================== Row 24 =====================
# becomes: Rosenticka / Rhodofomes roseus, A=111957843, E=658
$ws.Range("A24").Value = 111957843
$ws.Range("B24").Value = 89686
$ws.Range("E24").Value = 658
$ws.Range("F24").Value = "Rosenticka"
$ws.Range("G24").Value = "Rhodofomes roseus"
$ws.Range("H24").Value = "(Alb. & Schwein.) Kotl. & Pouzar"
Set-TextCell $ws.Range("I24") "6"
Set-TextCell $ws.Range("J24") "fruktkroppar"
$ws.Range("L24").ClearContents()
$ws.Range("M24").ClearContents()
$ws.Range("P24").Value = "Österås, Österås, Ång"
$ws.Range("Q24").Value = 609773.4136058968
$ws.Range("R24").Value = 7011992.49874373
$ws.Range("AC24").ClearContents()
Set-EmptyTextCell $ws.Range("AF24")

# ===================== Row 25 =====================
# becomes: Järpe / Tetrastes bonasia, A=111958182, E=102612
$ws.Range("A25").Value = 111958182
$ws.Range("B25").Value = 55611
$ws.Range("E25").Value = 102612
$ws.Range("F25").Value = "Järpe"
$ws.Range("G25").Value = "Tetrastes bonasia"
$ws.Range("H25").Value = "(Linnaeus, 1758)"
Set-TextCell $ws.Range("I25") "3"
$ws.Range("J25").ClearContents()
$ws.Range("L25").Value = "hona"
Set-EmptyTextCell $ws.Range("M25")
Set-TextCell $ws.Range("AC25") "1K"
$ws.Range("AF25").ClearContents()

# ===================== Row 27 =====================
# becomes: Rosenticka / Rhodofomes roseus, A=111957798 (loc: Österåsen)
$ws.Range("A27").Value = 111957798
$ws.Range("P27").Value = "Österåsen, Österås, Ång"
$ws.Range("Q27").Value = 609746.731343443
$ws.Range("R27").Value = 7011953.229753771
